# This script applies a re-shuffle of the observation records in rows 6-23:
# each row is updated in place with the field values (Id, Taxonsorteringsordning,
# Rodlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord, Start/Slutdatum)
# belonging to a different observation, per the source diff. Location/observer columns
# (P, T, U, V, W, AW, AX, ...) are left untouched since they are unchanged by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Startdatum/Slutdatum cells that are being rewritten as plain text
# (as in the original file) instead of letting Excel auto-convert the
# "yyyy-mm-dd" strings into date serial numbers.
$dateCells = @("Y8","AA8","Y11","AA11","Y14","AA14","Y15","AA15","Y16","AA16","Y17","AA17","Y18","AA18","Y20","AA20","Y22","AA22","Y23","AA23")
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 6
$ws.Range("A6").Value2 = 112475114
$ws.Range("B6").Value2 = 78746
$ws.Range("D6").Value2 = "LC"
$ws.Range("E6").Value2 = 6463
$ws.Range("F6").Value2 = "Bårdlav"
$ws.Range("G6").Value2 = "Nephroma parile"
$ws.Range("H6").Value2 = "(Ach.) Ach."
$ws.Range("Q6").Value2 = 537183
$ws.Range("R6").Value2 = 7202916

# Row 7
$ws.Range("A7").Value2 = 112475117
$ws.Range("B7").Value2 = 73834
$ws.Range("E7").Value2 = 6440
$ws.Range("F7").Value2 = "Vitgrynig nållav"
$ws.Range("G7").Value2 = "Chaenotheca subroscida"
$ws.Range("H7").Value2 = "(Eitner) Zahlbr."
$ws.Range("Q7").Value2 = 537176
$ws.Range("R7").Value2 = 7202935

# Row 8
$ws.Range("A8").Value2 = 112475125
$ws.Range("B8").Value2 = 89571
$ws.Range("D8").Value2 = "NT"
$ws.Range("E8").Value2 = 5432
$ws.Range("F8").Value2 = "Granticka"
$ws.Range("G8").Value2 = "Porodaedalea chrysoloma"
$ws.Range("H8").Value2 = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q8").Value2 = 537134
$ws.Range("R8").Value2 = 7202956
$ws.Range("Y8").Value2 = "2023-09-30"
$ws.Range("AA8").Value2 = "2023-09-30"

# Row 9
$ws.Range("A9").Value2 = 112475120
$ws.Range("B9").Value2 = 77650
$ws.Range("D9").Value2 = "NT"
$ws.Range("E9").Value2 = 6425
$ws.Range("F9").Value2 = "Garnlav"
$ws.Range("G9").Value2 = "Alectoria sarmentosa"
$ws.Range("Q9").Value2 = 537142
$ws.Range("R9").Value2 = 7202937

# Row 10
$ws.Range("A10").Value2 = 112475115
$ws.Range("Q10").Value2 = 537179
$ws.Range("R10").Value2 = 7202921

# Row 11
$ws.Range("A11").Value2 = 112475126
$ws.Range("B11").Value2 = 89553
$ws.Range("D11").Value2 = "NT"
$ws.Range("E11").Value2 = 1202
$ws.Range("F11").Value2 = "Ullticka"
$ws.Range("G11").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H11").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q11").Value2 = 537129
$ws.Range("R11").Value2 = 7202963
$ws.Range("Y11").Value2 = "2023-09-30"
$ws.Range("AA11").Value2 = "2023-09-30"

# Row 12
$ws.Range("A12").Value2 = 112475116
$ws.Range("B12").Value2 = 78740
$ws.Range("D12").Value2 = "LC"
$ws.Range("E12").Value2 = 6462
$ws.Range("F12").Value2 = "Stuplav"
$ws.Range("G12").Value2 = "Nephroma bellum"
$ws.Range("H12").Value2 = "(Spreng.) Tuck."
$ws.Range("R12").Value2 = 7202929

# Row 13
$ws.Range("A13").Value2 = 112475112
$ws.Range("B13").Value2 = 89571
$ws.Range("D13").Value2 = "NT"
$ws.Range("E13").Value2 = 5432
$ws.Range("F13").Value2 = "Granticka"
$ws.Range("G13").Value2 = "Porodaedalea chrysoloma"
$ws.Range("H13").Value2 = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q13").Value2 = 537180
$ws.Range("R13").Value2 = 7202916

# Row 14
$ws.Range("A14").Value2 = 112475119
$ws.Range("B14").Value2 = 89517
$ws.Range("E14").Value2 = 5447
$ws.Range("F14").Value2 = "Vedticka"
$ws.Range("G14").Value2 = "Fuscoporia viticola"
$ws.Range("H14").Value2 = "(Schwein.) Murrill"
$ws.Range("Q14").Value2 = 537157
$ws.Range("R14").Value2 = 7202946
$ws.Range("Y14").Value2 = "2023-09-30"
$ws.Range("AA14").Value2 = "2023-09-30"

# Row 15
$ws.Range("A15").Value2 = 112475046
$ws.Range("B15").Value2 = 97066
$ws.Range("D15").Value2 = "LC"
$ws.Range("E15").Value2 = 219880
$ws.Range("F15").Value2 = "Kransrams"
$ws.Range("G15").Value2 = "Polygonatum verticillatum"
$ws.Range("H15").Value2 = "(L.) All."
$ws.Range("Q15").Value2 = 537143
$ws.Range("R15").Value2 = 7202996
$ws.Range("Y15").Value2 = "2023-10-01"
$ws.Range("AA15").Value2 = "2023-10-01"

# Row 16
$ws.Range("A16").Value2 = 112475118
$ws.Range("B16").Value2 = 77650
$ws.Range("E16").Value2 = 6425
$ws.Range("F16").Value2 = "Garnlav"
$ws.Range("G16").Value2 = "Alectoria sarmentosa"
$ws.Range("H16").Value2 = "(Ach.) Ach."
$ws.Range("Q16").Value2 = 537145
$ws.Range("R16").Value2 = 7202941
$ws.Range("Y16").Value2 = "2023-09-30"
$ws.Range("AA16").Value2 = "2023-09-30"

# Row 17
$ws.Range("A17").Value2 = 112475048
$ws.Range("B17").Value2 = 89499
$ws.Range("D17").Value2 = "NT"
$ws.Range("E17").Value2 = 112
$ws.Range("F17").Value2 = "Stjärntagging"
$ws.Range("G17").Value2 = "Asterodon ferruginosus"
$ws.Range("H17").Value2 = "Pat."
$ws.Range("Q17").Value2 = 537193
$ws.Range("R17").Value2 = 7202876
$ws.Range("Y17").Value2 = "2023-09-30"
$ws.Range("AA17").Value2 = "2023-09-30"

# Row 18
$ws.Range("A18").Value2 = 112475044
$ws.Range("B18").Value2 = 78746
$ws.Range("D18").Value2 = "LC"
$ws.Range("E18").Value2 = 6463
$ws.Range("F18").Value2 = "Bårdlav"
$ws.Range("G18").Value2 = "Nephroma parile"
$ws.Range("H18").Value2 = "(Ach.) Ach."
$ws.Range("Q18").Value2 = 537164
$ws.Range("R18").Value2 = 7203017
$ws.Range("Y18").Value2 = "2023-10-01"
$ws.Range("AA18").Value2 = "2023-10-01"

# Row 20
$ws.Range("A20").Value2 = 112475042
$ws.Range("B20").Value2 = 93324
$ws.Range("D20").Value2 = "VU"
$ws.Range("E20").Value2 = 1079
$ws.Range("F20").Value2 = "Aspfjädermossa"
$ws.Range("G20").Value2 = "Neckera pennata"
$ws.Range("H20").Value2 = "Hedw."
$ws.Range("Q20").Value2 = 537165
$ws.Range("R20").Value2 = 7203022
$ws.Range("Y20").Value2 = "2023-10-01"
$ws.Range("AA20").Value2 = "2023-10-01"

# Row 22
$ws.Range("A22").Value2 = 112475047
$ws.Range("B22").Value2 = 89567
$ws.Range("E22").Value2 = 1204
$ws.Range("F22").Value2 = "Gränsticka"
$ws.Range("G22").Value2 = "Phellopilus nigrolimitatus"
$ws.Range("H22").Value2 = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q22").Value2 = 537125
$ws.Range("R22").Value2 = 7202974
$ws.Range("Y22").Value2 = "2023-10-01"
$ws.Range("AA22").Value2 = "2023-10-01"

# Row 23
$ws.Range("A23").Value2 = 112475043
$ws.Range("B23").Value2 = 78740
$ws.Range("E23").Value2 = 6462
$ws.Range("F23").Value2 = "Stuplav"
$ws.Range("G23").Value2 = "Nephroma bellum"
$ws.Range("H23").Value2 = "(Spreng.) Tuck."
$ws.Range("Q23").Value2 = 537164
$ws.Range("R23").Value2 = 7203021
$ws.Range("Y23").Value2 = "2023-10-01"
$ws.Range("AA23").Value2 = "2023-10-01"
